$p = $ppt.ActivePresentation

# --- 1. Table on slide 5 (the B1/B2 financial-documents table): apply the
#        new gallery table style. PowerPoint COM requires ApplyStyle(), a
#        plain property assignment on Table.Style is rejected.
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{7277B8A9-E766-4EEB-A514-B2D0266B25C9}")

# --- 2. Switch the deck's design colours from the "Red Violet" (Integral)
#        palette over to the standard "Office" palette - i.e. the master's
#        theme colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#        COM exposes the DrawingML 12-slot colour scheme through
#        Theme.ThemeColorScheme; each slot's .RGB takes a COLORREF
#        (0x00BBGGRR), so we convert from the RRGGBB hex values below.
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function ToComRgb([int]$rrggbb) {
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = ToComRgb($officeColors[$i])
}
